$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 500
$ws.Range("I2").Value = 450
$ws.Range("J2").Value = 550
$ws.Range("K2").Value = 450
$ws.Range("L2").Value = 550
$ws.Range("M2").Value = -337
$ws.Range("N2").Value = -776

$ws.Range("H6").Value = 254.4
$ws.Range("I6").Value = 236.75
$ws.Range("J6").Value = 325
$ws.Range("K6").Value = 710.25
$ws.Range("L6").Value = 975
$ws.Range("M6").Value = -598.25
$ws.Range("N6").Value = -1199

$ws.Range("H29").Value = 733.2222
$ws.Range("I29").Value = 199.875
$ws.Range("J29").Value = 5000
$ws.Range("K29").Value = 599.625
$ws.Range("L29").Value = 15000
$ws.Range("M29").Value = -318.625
$ws.Range("N29").Value = -15562

$ws.Range("H38").Value = 644.3333
$ws.Range("I38").Value = 53.2
$ws.Range("J38").Value = 3600
$ws.Range("K38").Value = 159.6
$ws.Range("L38").Value = 10800
$ws.Range("M38").Value = 212.4
$ws.Range("N38").Value = -11544

$ws.Range("H40").Value = 7649.9
$ws.Range("I40").Value = 4875
$ws.Range("J40").Value = 8343.625
$ws.Range("K40").Value = 4875
$ws.Range("L40").Value = 8343.625
$ws.Range("M40").Value = -4700
$ws.Range("N40").Value = -8693.625

$ws.Range("H43").Value = 1399.5
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 1399.5
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 1399.5
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -1537.5

$ws.Range("H55").Value = 259
$ws.Range("I55").Value = 198.33333
$ws.Range("J55").Value = 350
$ws.Range("K55").Value = 198.33333
$ws.Range("L55").Value = 350
$ws.Range("M55").Value = 15.66667000000001
$ws.Range("N55").Value = -778

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()

$ws.Range("H135").Value = 2031
$ws.Range("I135").Value = 1625.1666
$ws.Range("J135").Value = 3005
$ws.Range("K135").Value = 14626.4994
$ws.Range("L135").Value = 27045
$ws.Range("M135").Value = -12091.4994
$ws.Range("N135").Value = -32115

$ws.Range("H138").Value = 10098
$ws.Range("I138").Value = 5197
$ws.Range("J138").Value = 14999
$ws.Range("K138").Value = 15591
$ws.Range("L138").Value = 44997
$ws.Range("M138").Value = -10451
$ws.Range("N138").Value = -55277

$ws.Range("H141").Value = 933.4
$ws.Range("I141").Value = 933.4
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 2800.2
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 2379.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 5666.6665
$ws.Range("I15").Value = 1000
$ws.Range("J15").Value = 15000
$ws.Range("K15").Value = 1000
$ws.Range("L15").Value = 15000
$ws.Range("M15").Value = -650
$ws.Range("N15").Value = -15700

$ws.Range("H61").Value = 4972.5
$ws.Range("I61").Value = 4972.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 4972.5
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -4760.5

$ws.Range("H74").Value = 1100.1
$ws.Range("I74").Value = 1100.1
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1100.1
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -226.0999999999999

$ws.Range("H77").Value = 1100.1
$ws.Range("I77").Value = 1100.1
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 5500.5
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -1132.5

$ws.Range("H122").Value = 1749
$ws.Range("I122").Value = 1749
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5247
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2797
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 4951.7
$ws.Range("I132").Value = 3190
$ws.Range("J132").Value = 11998.5
$ws.Range("K132").Value = 9570
$ws.Range("L132").Value = 35995.5
$ws.Range("M132").Value = -7040
$ws.Range("N132").Value = -41055.5

$ws.Range("H136").Value = 4972.5
$ws.Range("I136").Value = 4972.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 14917.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -12367.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 14849.4
$ws.Range("I82").Value = 14849.4
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 14849.4
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -14466.4

$ws.Range("H85").Value = 14849.4
$ws.Range("I85").Value = 14849.4
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 14849.4
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -13523.4

$ws.Range("H86").Value = 4566.923
$ws.Range("I86").Value = 2329.6365
$ws.Range("J86").Value = 6207.6
$ws.Range("K86").Value = 2329.6365
$ws.Range("L86").Value = 6207.6
$ws.Range("M86").Value = -1206.6365
$ws.Range("N86").Value = -8453.6

$ws.Range("H89").Value = 4566.923
$ws.Range("I89").Value = 2329.6365
$ws.Range("J89").Value = 6207.6
$ws.Range("K89").Value = 11648.1825
$ws.Range("L89").Value = 31038
$ws.Range("M89").Value = -6032.182500000001
$ws.Range("N89").Value = -42270

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 250
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 250
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 250
$ws.Range("N22").Value = -950

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H107").Value = 741.86664
$ws.Range("I107").Value = 684.25
$ws.Range("J107").Value = 807.7143
$ws.Range("K107").Value = 684.25
$ws.Range("L107").Value = 807.7143
$ws.Range("M107").Value = 1235.75
$ws.Range("N107").Value = -4647.7143

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 85
$ws.Range("I2").Value = 90
$ws.Range("J2").Value = 81.666664
$ws.Range("K2").Value = 540
$ws.Range("L2").Value = 489.999984
$ws.Range("M2").Value = -427
$ws.Range("N2").Value = -715.999984

$ws.Range("H23").Value = 667
$ws.Range("I23").Value = 157
$ws.Range("J23").Value = 1177
$ws.Range("K23").Value = 471
$ws.Range("L23").Value = 3531
$ws.Range("M23").Value = -236
$ws.Range("N23").Value = -4001

$ws.Range("H34").Value = 4000
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 4000
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 12000
$ws.Range("N34").Value = -12168

$ws.Range("H39").Value = 905.44446
$ws.Range("I39").Value = 850
$ws.Range("J39").Value = 912.375
$ws.Range("K39").Value = 2550
$ws.Range("L39").Value = 2737.125
$ws.Range("M39").Value = -2256
$ws.Range("N39").Value = -3325.125

$ws.Range("H118").Value = 500
$ws.Range("I118").Value = 500
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 1500
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -257

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 2504000
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 2504000
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 2504000
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -2504224

$ws.Range("H8").Value = 2504000
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 2504000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 2504000
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -2504278

$ws.Range("H43").Value = 10166.667
$ws.Range("I43").Value = 500
$ws.Range("J43").Value = 15000
$ws.Range("K43").Value = 500
$ws.Range("L43").Value = 15000
$ws.Range("M43").Value = -349
$ws.Range("N43").Value = -15302

$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H122").Value = 13756.667
$ws.Range("I122").Value = 9517.467000000001
$ws.Range("J122").Value = 34952.668
$ws.Range("K122").Value = 28552.401
$ws.Range("L122").Value = 104858.004
$ws.Range("M122").Value = -26102.401
$ws.Range("N122").Value = -109758.004

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 694.7
$ws.Range("I16").Value = 660.7778
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 660.7778
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -490.7778
$ws.Range("N16").Value = -1340

$ws.Range("H46").Value = 1614.1936
$ws.Range("I46").Value = 1062.1333
$ws.Range("J46").Value = 2131.75
$ws.Range("K46").Value = 1062.1333
$ws.Range("L46").Value = 2131.75
$ws.Range("M46").Value = -874.1333
$ws.Range("N46").Value = -2507.75

$ws.Range("H68").Value = 2871.1428
$ws.Range("I68").Value = 2959.6
$ws.Range("J68").Value = 2650
$ws.Range("K68").Value = 2959.6
$ws.Range("L68").Value = 2650
$ws.Range("M68").Value = -2210.6
$ws.Range("N68").Value = -4148

$ws.Range("H71").Value = 2871.1428
$ws.Range("I71").Value = 2959.6
$ws.Range("J71").Value = 2650
$ws.Range("K71").Value = 14798
$ws.Range("L71").Value = 13250
$ws.Range("M71").Value = -11054
$ws.Range("N71").Value = -20738

$ws.Range("H82").Value = 2229.7
$ws.Range("I82").Value = 2310.7778
$ws.Range("J82").Value = 1500
$ws.Range("K82").Value = 2310.7778
$ws.Range("L82").Value = 1500
$ws.Range("M82").Value = -1949.7778
$ws.Range("N82").Value = -2222

$ws.Range("H85").Value = 2229.7
$ws.Range("I85").Value = 2310.7778
$ws.Range("J85").Value = 1500
$ws.Range("K85").Value = 2310.7778
$ws.Range("L85").Value = 1500
$ws.Range("M85").Value = -1062.7778
$ws.Range("N85").Value = -3996

$ws.Range("H132").Value = 38834
$ws.Range("I132").Value = 54501.332
$ws.Range("J132").Value = 23166.666
$ws.Range("K132").Value = 163503.996
$ws.Range("L132").Value = 69499.99800000001
$ws.Range("M132").Value = -160973.996
$ws.Range("N132").Value = -74559.99800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()

$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()

$ws.Range("H62").Value = 4200
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 4200
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 4200
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -5448

$ws.Range("H65").Value = 4200
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 4200
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 21000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -27240

$ws.Range("H107").Value = 1272.25
$ws.Range("I107").Value = 738.4286
$ws.Range("J107").Value = 2019.6
$ws.Range("K107").Value = 2215.2858
$ws.Range("L107").Value = 6058.799999999999
$ws.Range("M107").Value = -295.2857999999997
$ws.Range("N107").Value = -9898.799999999999

$ws.Range("H113").Value = 615.6667
$ws.Range("I113").Value = 624.625
$ws.Range("J113").Value = 597.75
$ws.Range("K113").Value = 1873.875
$ws.Range("L113").Value = 1793.25
$ws.Range("M113").Value = 624.625
$ws.Range("N113").Value = -6133.25
